$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7").Value = "This row has a`nmulti-line value"
$ws.Range("A7").WrapText = $true
$ws.Range("A7").Select()
